# Auto-generated edit script: applies the ifrs list correction for 유엔젤.xlsx
# Rows 2-6: replace financial figures in columns D:AJ (some cells in rows 4-6 are cleared: AG/AH)
# Rows 7-9: clear all financial figures in columns D:AJ, keeping only columns A-C
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$r2_D = New-Object 'object[,]' 1,33
$r2_D[0,0] = 354
$r2_D[0,1] = -42
$r2_D[0,2] = -42
$r2_D[0,3] = -68
$r2_D[0,4] = -96
$r2_D[0,5] = -97
$r2_D[0,6] = 0
$r2_D[0,7] = 602
$r2_D[0,8] = 72
$r2_D[0,9] = 530
$r2_D[0,10] = 526
$r2_D[0,11] = 4
$r2_D[0,12] = 66
$r2_D[0,13] = -30
$r2_D[0,14] = 72
$r2_D[0,15] = -22
$r2_D[0,16] = 4
$r2_D[0,17] = -34
$r2_D[0,18] = 30
$r2_D[0,19] = -11.83
$r2_D[0,20] = -27.24
$r2_D[0,21] = -16.54
$r2_D[0,22] = -14.51
$r2_D[0,23] = 13.63
$r2_D[0,24] = 844.22
$r2_D[0,25] = -732
$r2_D[0,26] = -6.43
$r2_D[0,27] = 4578
$r2_D[0,28] = 1.03
$r2_D[0,29] = 100
$r2_D[0,30] = 2.13
$r2_D[0,31] = -11.89
$r2_D[0,32] = 13195454
$ws.Range("D2:AJ2").Value = $r2_D

# Row 3
$r3_D = New-Object 'object[,]' 1,33
$r3_D[0,0] = 412
$r3_D[0,1] = 6
$r3_D[0,2] = 6
$r3_D[0,3] = 19
$r3_D[0,4] = 11
$r3_D[0,5] = 15
$r3_D[0,6] = -4
$r3_D[0,7] = 621
$r3_D[0,8] = 94
$r3_D[0,9] = 527
$r3_D[0,10] = 526
$r3_D[0,11] = 1
$r3_D[0,12] = 66
$r3_D[0,13] = 8
$r3_D[0,14] = 36
$r3_D[0,15] = -12
$r3_D[0,16] = 4
$r3_D[0,17] = 4
$r3_D[0,18] = 30
$r3_D[0,19] = 1.56
$r3_D[0,20] = 2.68
$r3_D[0,21] = 2.79
$r3_D[0,22] = 1.81
$r3_D[0,23] = 17.84
$r3_D[0,24] = 844.77
$r3_D[0,25] = 111
$r3_D[0,26] = 51.19
$r3_D[0,27] = 4580
$r3_D[0,28] = 1.24
$r3_D[0,29] = 100
$r3_D[0,30] = 1.76
$r3_D[0,31] = 78.43000000000001
$r3_D[0,32] = 13195454
$ws.Range("D3:AJ3").Value = $r3_D

# Row 4
$r4_D = New-Object 'object[,]' 1,29
$r4_D[0,0] = 303
$r4_D[0,1] = -47
$r4_D[0,2] = -47
$r4_D[0,3] = -47
$r4_D[0,4] = -47
$r4_D[0,5] = -43
$r4_D[0,6] = -5
$r4_D[0,7] = 583
$r4_D[0,8] = 120
$r4_D[0,9] = 463
$r4_D[0,10] = 467
$r4_D[0,11] = -4
$r4_D[0,12] = 66
$r4_D[0,13] = -36
$r4_D[0,14] = -17
$r4_D[0,15] = 8
$r4_D[0,16] = 4
$r4_D[0,17] = -40
$r4_D[0,18] = 50
$r4_D[0,19] = -15.59
$r4_D[0,20] = -15.65
$r4_D[0,21] = -8.630000000000001
$r4_D[0,22] = -7.89
$r4_D[0,23] = 25.98
$r4_D[0,24] = 762.51
$r4_D[0,25] = -325
$r4_D[0,26] = -15.65
$r4_D[0,27] = 4066
$r4_D[0,28] = 1.25
$ws.Range("D4:AF4").Value = $r4_D
$r4_AI = New-Object 'object[,]' 1,2
$r4_AI[0,0] = 0
$r4_AI[0,1] = 13195454
$ws.Range("AI4:AJ4").Value = $r4_AI
$ws.Range("AG4:AH4").ClearContents()

# Row 5
$r5_D = New-Object 'object[,]' 1,29
$r5_D[0,0] = 319
$r5_D[0,1] = -10
$r5_D[0,2] = -10
$r5_D[0,3] = -39
$r5_D[0,4] = -66
$r5_D[0,5] = -64
$r5_D[0,6] = -2
$r5_D[0,7] = 512
$r5_D[0,8] = 114
$r5_D[0,9] = 398
$r5_D[0,10] = 405
$r5_D[0,11] = -7
$r5_D[0,12] = 66
$r5_D[0,13] = -20
$r5_D[0,14] = 31
$r5_D[0,15] = 0
$r5_D[0,16] = 4
$r5_D[0,17] = -24
$r5_D[0,18] = 50
$r5_D[0,19] = -3.08
$r5_D[0,20] = -20.75
$r5_D[0,21] = -14.67
$r5_D[0,22] = -12.1
$r5_D[0,23] = 28.65
$r5_D[0,24] = 665.52
$r5_D[0,25] = -485
$r5_D[0,26] = -9.07
$r5_D[0,27] = 3531
$r5_D[0,28] = 1.24
$ws.Range("D5:AF5").Value = $r5_D
$r5_AI = New-Object 'object[,]' 1,2
$r5_AI[0,0] = 0
$r5_AI[0,1] = 13195454
$ws.Range("AI5:AJ5").Value = $r5_AI
$ws.Range("AG5:AH5").ClearContents()

# Row 6
$r6_D = New-Object 'object[,]' 1,6
$r6_D[0,0] = 335
$r6_D[0,1] = 30
$r6_D[0,2] = 30
$r6_D[0,3] = 21
$r6_D[0,4] = 18
$r6_D[0,5] = 18
$ws.Range("D6:I6").Value = $r6_D
$r6_K = New-Object 'object[,]' 1,4
$r6_K[0,0] = 512
$r6_K[0,1] = 102
$r6_K[0,2] = 410
$r6_K[0,3] = 416
$ws.Range("K6:N6").Value = $r6_K
$r6_P = New-Object 'object[,]' 1,17
$r6_P[0,0] = 66
$r6_P[0,1] = 31
$r6_P[0,2] = 86
$r6_P[0,3] = -20
$r6_P[0,4] = 4
$r6_P[0,5] = 27
$r6_P[0,6] = 30
$r6_P[0,7] = 8.99
$r6_P[0,8] = 5.49
$r6_P[0,9] = 4.35
$r6_P[0,10] = 3.59
$r6_P[0,11] = 24.91
$r6_P[0,12] = 738.76
$r6_P[0,13] = 135
$r6_P[0,14] = 32.31
$r6_P[0,15] = 3619
$r6_P[0,16] = 1.21
$ws.Range("P6:AF6").Value = $r6_P
$r6_AI = New-Object 'object[,]' 1,2
$r6_AI[0,0] = 0
$r6_AI[0,1] = 13195454
$ws.Range("AI6:AJ6").Value = $r6_AI
$ws.Range("AG6:AH6").ClearContents()

# Rows 7-9: clear all data columns (D:AJ), keeping A-C (rank/period label/year label)
$ws.Range("D7:AJ9").ClearContents()
